$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3244.375
$ws.Range("I100").Value = 2876.25
$ws.Range("J100").Value = 3612.5
$ws.Range("K100").Value = 2876.25
$ws.Range("L100").Value = 3612.5
$ws.Range("M100").Value = -2335.25
$ws.Range("N100").Value = -4694.5

$ws.Range("H111").Value = 2684.2727
$ws.Range("I111").Value = 1704.5
$ws.Range("J111").Value = 3860
$ws.Range("K111").Value = 5113.5
$ws.Range("L111").Value = 11580
$ws.Range("M111").Value = -2046.5
$ws.Range("N111").Value = -17714

$ws.Range("H132").Value = 3263.1072
$ws.Range("I132").Value = 3324.125
$ws.Range("J132").Value = 2897
$ws.Range("K132").Value = 9972.375
$ws.Range("L132").Value = 8691
$ws.Range("M132").Value = -7442.375
$ws.Range("N132").Value = -13751

$ws.Range("H137").Value = 1870.8422
$ws.Range("I137").Value = 1886
$ws.Range("J137").Value = 1850
$ws.Range("K137").Value = 5658
$ws.Range("L137").Value = 5550
$ws.Range("M137").Value = -3108
$ws.Range("N137").Value = -10650

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1701

$ws.Range("H32").Value = 3684.7812
$ws.Range("I32").Value = 3319.111
$ws.Range("J32").Value = 5659.4
$ws.Range("K32").Value = 3319.111
$ws.Range("L32").Value = 5659.4
$ws.Range("M32").Value = -3032.111
$ws.Range("N32").Value = -6233.4

$ws.Range("H61").Value = 2235.8333
$ws.Range("I61").Value = 1240.5416
$ws.Range("J61").Value = 4226.4165
$ws.Range("K61").Value = 1240.5416
$ws.Range("L61").Value = 4226.4165
$ws.Range("M61").Value = -1028.5416
$ws.Range("N61").Value = -4650.4165

$ws.Range("H74").Value = 2836.3809
$ws.Range("I74").Value = 3007.1714
$ws.Range("J74").Value = 1982.4286
$ws.Range("K74").Value = 3007.1714
$ws.Range("L74").Value = 1982.4286
$ws.Range("M74").Value = -2133.1714
$ws.Range("N74").Value = -3730.4286

$ws.Range("H77").Value = 2836.3809
$ws.Range("I77").Value = 3007.1714
$ws.Range("J77").Value = 1982.4286
$ws.Range("K77").Value = 15035.857
$ws.Range("L77").Value = 9912.143
$ws.Range("M77").Value = -10667.857
$ws.Range("N77").Value = -18648.143

$ws.Range("H101").Value = 33840.6
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 33840.6
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 33840.6
$ws.Range("N101").Value = -40330.6

$ws.Range("H102").Value = 2143.1428
$ws.Range("I102").Value = 2200.4
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2200.4
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -578.4000000000001
$ws.Range("N102").Value = -5244

$ws.Range("H136").Value = 2235.8333
$ws.Range("I136").Value = 1240.5416
$ws.Range("J136").Value = 4226.4165
$ws.Range("K136").Value = 3721.6248
$ws.Range("L136").Value = 12679.2495
$ws.Range("M136").Value = -1171.6248
$ws.Range("N136").Value = -17779.2495

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1727.56
$ws.Range("I94").Value = 1389
$ws.Range("J94").Value = 3505
$ws.Range("K94").Value = 1389
$ws.Range("L94").Value = 3505
$ws.Range("M94").Value = -938
$ws.Range("N94").Value = -4407

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1566.1538
$ws.Range("I5").Value = 1450.909
$ws.Range("J5").Value = 2200
$ws.Range("K5").Value = 4352.727000000001
$ws.Range("L5").Value = 6600
$ws.Range("M5").Value = -4240.727000000001
$ws.Range("N5").Value = -6824

$ws.Range("H33").Value = 229.8
$ws.Range("I33").Value = 183
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 1098
$ws.Range("L33").Value = 1800
$ws.Range("M33").Value = -815
$ws.Range("N33").Value = -2366

$ws.Range("H122").Value = 704.6667
$ws.Range("I122").Value = 416.66666
$ws.Range("J122").Value = 992.6667
$ws.Range("K122").Value = 3749.99994
$ws.Range("L122").Value = 8934.0003
$ws.Range("M122").Value = -1299.99994
$ws.Range("N122").Value = -13834.0003

$ws.Range("H131").Value = 813.27
$ws.Range("I131").Value = 750
$ws.Range("J131").Value = 815.90625
$ws.Range("K131").Value = 2250
$ws.Range("L131").Value = 2447.71875
$ws.Range("M131").Value = 2790
$ws.Range("N131").Value = -12527.71875

$ws.Range("H135").Value = 1566.1538
$ws.Range("I135").Value = 1450.909
$ws.Range("J135").Value = 2200
$ws.Range("K135").Value = 13058.181
$ws.Range("L135").Value = 19800
$ws.Range("M135").Value = -10523.181
$ws.Range("N135").Value = -24870

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3649.75
$ws.Range("I80").Value = 3200.6667
$ws.Range("J80").Value = 4227.143
$ws.Range("K80").Value = 3200.6667
$ws.Range("L80").Value = 4227.143
$ws.Range("M80").Value = -2202.6667
$ws.Range("N80").Value = -6223.143

$ws.Range("H83").Value = 3649.75
$ws.Range("I83").Value = 3200.6667
$ws.Range("J83").Value = 4227.143
$ws.Range("K83").Value = 16003.3335
$ws.Range("L83").Value = 21135.715
$ws.Range("M83").Value = -11011.3335
$ws.Range("N83").Value = -31119.715

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2430
$ws.Range("I22").Value = 1905.1
$ws.Range("J22").Value = 3479.8
$ws.Range("K22").Value = 1905.1
$ws.Range("L22").Value = 3479.8
$ws.Range("M22").Value = -1610.1
$ws.Range("N22").Value = -4069.8

$ws.Range("H27").Value = 2430
$ws.Range("I27").Value = 1905.1
$ws.Range("J27").Value = 3479.8
$ws.Range("K27").Value = 1905.1
$ws.Range("L27").Value = 3479.8
$ws.Range("M27").Value = -1798.1
$ws.Range("N27").Value = -3693.8

$ws.Range("H40").Value = 3740.2
$ws.Range("I40").Value = 3030
$ws.Range("J40").Value = 4095.3
$ws.Range("K40").Value = 3030
$ws.Range("L40").Value = 4095.3
$ws.Range("M40").Value = -2894
$ws.Range("N40").Value = -4367.3

$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -3376

$ws.Range("H68").Value = 2584.8572
$ws.Range("I68").Value = 2531.5
$ws.Range("J68").Value = 2624.875
$ws.Range("K68").Value = 2531.5
$ws.Range("L68").Value = 2624.875
$ws.Range("M68").Value = -1782.5
$ws.Range("N68").Value = -4122.875

$ws.Range("H71").Value = 2584.8572
$ws.Range("I71").Value = 2531.5
$ws.Range("J71").Value = 2624.875
$ws.Range("K71").Value = 12657.5
$ws.Range("L71").Value = 13124.375
$ws.Range("M71").Value = -8913.5
$ws.Range("N71").Value = -20612.375

$ws.Range("H82").Value = 4239
$ws.Range("I82").Value = 5000
$ws.Range("J82").Value = 3097.5
$ws.Range("K82").Value = 5000
$ws.Range("L82").Value = 3097.5
$ws.Range("M82").Value = -4639
$ws.Range("N82").Value = -3819.5

$ws.Range("H85").Value = 4239
$ws.Range("I85").Value = 5000
$ws.Range("J85").Value = 3097.5
$ws.Range("K85").Value = 5000
$ws.Range("L85").Value = 3097.5
$ws.Range("M85").Value = -3752
$ws.Range("N85").Value = -5593.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H62").Value = 2999
$ws.Range("I62").Value = 2999
$ws.Range("J62").Value = 2999
$ws.Range("K62").Value = 2999
$ws.Range("L62").Value = 2999
$ws.Range("M62").Value = -2375
$ws.Range("N62").Value = -4247

$ws.Range("H65").Value = 2999
$ws.Range("I65").Value = 2999
$ws.Range("J65").Value = 2999
$ws.Range("K65").Value = 14995
$ws.Range("L65").Value = 14995
$ws.Range("M65").Value = -11875
$ws.Range("N65").Value = -21235

$ws.Range("H100").Value = 866.44446
$ws.Range("I100").Value = 649.75
$ws.Range("J100").Value = 1039.8
$ws.Range("K100").Value = 1299.5
$ws.Range("L100").Value = 2079.6
$ws.Range("M100").Value = -758.5
$ws.Range("N100").Value = -3161.6

$ws.Range("H116").Value = 37495
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 37495
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 37495
$ws.Range("N116").Value = -46673

$ws.Range("H117").Value = 40000
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 40000
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 40000
$ws.Range("N117").Value = -49178

$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H120").Value = 34420
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 34420
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 34420
$ws.Range("N120").Value = -44096

$ws.Range("H136").Value = 18183090
$ws.Range("I136").Value = 27028138
$ws.Range("J136").Value = 1603.1111
$ws.Range("K136").Value = 81084414
$ws.Range("L136").Value = 4809.3333
$ws.Range("M136").Value = -81081864
$ws.Range("N136").Value = -9909.3333
